$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume% (E) cells hold numeric-looking strings (e.g.
# "238.33", "30.394.32", "  -1.22%  ") that must stay literal TEXT, just
# like the original inlineStr cells -- not get auto-coerced into real
# numbers by Excel smart-entry parsing. Force text format first.
$forceTextCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6",
    "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11",
    "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16",
    "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21",
    "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26",
    "D27", "E27", "D28", "E28", "D29", "E29", "E30", "D31", "E31", "D32",
    "E32", "D33", "E33", "D34", "E34", "D35", "D36", "E36", "D37", "E37",
    "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43",
    "E43", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48",
    "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($c in $forceTextCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.394.32"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "1.891.67"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "238.33"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4820"
$ws.Range("E7").Value = "  -2.13%  "
$ws.Range("D8").Value = "0.2903"
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("D9").Value = "0.06606"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").Value = "1.895.38"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").Value = "16.94"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").Value = "0.07452"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "5.177"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "87.67"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").Value = "0.6642"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "30.368.70"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "13.45"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000007773"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "5.444"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "2.140.22"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "194.10"
$ws.Range("E23").Value = "  -4.08%  "
$ws.Range("D24").Value = "6.186"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").Value = "9.421"
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").Value = "163.54"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("D27").Value = "18.26"
$ws.Range("E27").Value = "  -3.62%  "
$ws.Range("D28").Value = "1.953"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").Value = "1.450"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").Value = "0.09143"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").Value = "4.065"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").Value = "0.05104"
$ws.Range("E33").Value = "  -3.89%  "
$ws.Range("D34").Value = "1.151"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "0.7334"
$ws.Range("D36").Value = "2.710"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").Value = "0.01817"
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("D39").Value = "0.9186"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").Value = "2.085"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "106.62"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.900"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").Value = "0.4336"
$ws.Range("E43").Value = "  -4.29%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "7.673"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").Value = "0.1334"
$ws.Range("E46").Value = "  -5.25%  "
$ws.Range("D47").Value = "1.563"
$ws.Range("E47").Value = "  +6.96%  "
$ws.Range("D48").Value = "64.92"
$ws.Range("E48").Value = "  -11.33%  "
$ws.Range("D49").Value = "8.981"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "34.10"
$ws.Range("E50").Value = "  -5.58%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05769"
$ws.Range("E51").Value = "  -2.90%  "

# Restore default (general) style on the force-text cells so only the
# value changed -- no lingering number-format override remains on them.
foreach ($c in $forceTextCells) {
    $ws.Range($c).Style = "Normal"
}
